$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 23, shifting existing rows 23-28 down to 24-29
$ws.Rows.Item(23).Insert()

# Fill the new row 23 with data
$ws.Range("A23").Value = 11
$ws.Range("B23").Value = "Vega Monumental Concepción"
$ws.Range("C23").Value = "Bíobío"
$ws.Range("D23").Value = 44510
$ws.Range("E23").Value = 8
$ws.Range("F23").Value = 100112031
$ws.Range("G23").Value = "Poroto verde"
$ws.Range("H23").Value = "Magnum"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 35000
$ws.Range("L23").Value = 36000
$ws.Range("M23").Value = 35500
$ws.Range("N23").Value = "$/malla 25 kilos"
$ws.Range("O23").Value = "Perú"
$ws.Range("P23").Value = 1420
$ws.Range("Q23").Value = 25
$ws.Range("R23").Value = "Hortaliza"
